$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three numeric cell values in row 1
$ws.Range("A1").Value = 148.97112728234609
$ws.Range("B1").Value = 4.5794831703104277
$ws.Range("C1").Value = 1.0074363992172213

# Narrow columns B and C by one character width each
# (B: 11.7109375 -> 10.7109375, C: 12.7109375 -> 11.7109375).
# ColumnWidth is stored/rounded internally to whole-pixel steps, so the
# assigned values are chosen to land the resulting (pixel-quantized)
# column width as close as possible to the target width above.
$ws.Columns.Item(2).ColumnWidth = 9.85
$ws.Columns.Item(3).ColumnWidth = 10.85
